$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Top Gainers")

# Shift rows 63-65 (ORIENTTECH, ICRA, SALASAR) down into rows 64-66,
# and place the updated NPST data into row 63.

$ws.Range("B63").Value = "NPST"
$ws.Range("C63").Value = 3.8509
$ws.Range("D63").Value = -2.0059
$ws.Range("E63").Value = -3.5057

$ws.Range("B64").Value = "ORIENTTECH"
$ws.Range("C64").Value = 3.827
$ws.Range("D64").Value = 0.5247000000000001
$ws.Range("E64").Value = 32.6784

$ws.Range("B65").Value = "ICRA"
$ws.Range("C65").Value = 3.7985
$ws.Range("D65").Value = 4.4793
$ws.Range("E65").Value = 2.8828

$ws.Range("B66").Value = "SALASAR"
$ws.Range("C66").Value = 3.7935
$ws.Range("D66").Value = 4.7872
$ws.Range("E66").Value = 11.0485
